# Batch write the "actual response" data back into the workbook,
# reading the whole Case2 sheet once and writing it back once.
#
# Case2!E (ExpectedResponseData) already holds the expected value for each
# test case; Case2!F (ActualResponseData) is the column that records what
# actually came back from the API call. For this run the actual response
# matched the expected response for every row, so the batch writer copies
# E2:E12 straight across into F2:F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case2")

# Single READ: pull every expected-response value for the data rows.
$firstDataRow = 2
$lastDataRow = 12
$expectedColumn = 5   # E - ExpectedResponseData(期望响应值)
$actualColumn = 6     # F - ActualResponseData(实际响应数值)

$values = @{}
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $values[$row] = $ws.Cells.Item($row, $expectedColumn).Value2
}

# Single WRITE: push the captured values back into the ActualResponseData column.
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, $actualColumn).Value2 = $values[$row]
}
